$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("P1 - Anastasia")
$ws2 = $wb.Worksheets.Item("P2 - Jan Willem")
$ws3 = $wb.Worksheets.Item("P3 - Ivar")
$ws4 = $wb.Worksheets.Item("P4 - Roy")
$ws5 = $wb.Worksheets.Item("P5 - Jasper")

# ---------------------------------------------------------------------------
# New logbook rows.
# Cells that introduce brand-new shared strings are written in the exact
# order needed so the shared-strings table ends up with the new entries in
# the same sequence as the target workbook (212..217).
# ---------------------------------------------------------------------------

# P2 - Jan Willem: row 57
$ws2.Range("A57").Value = "Bezig met de code "
$ws2.Range("B57").Value = 44167
$ws2.Range("B56").Copy()
$ws2.Range("B57").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws2.Range("C57").Value = 120
$ws2.Range("D57").Value = "Geprobeerd de customer service te maken, niet gelukt vraag morgen hulp van anderen"

# P5 - Jasper: row 71 (D cell written now to fix shared-string order)
$ws5.Range("D71").Value = "Start gemaakt met securityOpdracht.  Demo powerpoint gemaakt & uitlogknop gemaakt"

# P2 - Jan Willem: row 58
$ws2.Range("A58").Value = "KBS op locatie"
$ws2.Range("B58").Value = 44168
$ws2.Range("B56").Copy()
$ws2.Range("B58").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws2.Range("C58").Value = 180
$ws2.Range("D58").Value = "Start gemaakt met securityOpdracht. stijlen knoppen en versturen mail voor klantservice, powerpoint voor de demo gemaakt."

# P1 - Anastasia: row 54
$ws1.Range("D54").Value = "Securityverslag"
$ws1.Range("A54").Value = "KBS op locatie via teams"
$ws1.Range("B54").Value = 44168
$ws1.Range("B53").Copy()
$ws1.Range("B54").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws1.Range("C54").Value = 90

# P2 - Jan Willem: row 59
$ws2.Range("A59").Value = "Mail klantenservice"
$ws2.Range("B59").Value = 44168
$ws2.Range("B56").Copy()
$ws2.Range("B59").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws2.Range("C59").Value = 15
$ws2.Range("D59").Value = "sendMail.php meermaals geïncluded. Meermaalse verwijderd"

# P5 - Jasper: row 71 remaining cells
$ws5.Range("A71").Value = "KBS op locatie"
$ws5.Range("B71").Value = 44168
$ws5.Range("B70").Copy()
$ws5.Range("B71").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws5.Range("C71").Value = 180

# ---------------------------------------------------------------------------
# P5 - Jasper was missing the "Totaal" hour formula in B6; add it now so it
# matches the other sheets.
# ---------------------------------------------------------------------------
$ws5.Range("B6").Formula = "=SUM(C10:C186)/60"

# ---------------------------------------------------------------------------
# Grow the bound tables (ListObjects) so the new rows are included.
# ---------------------------------------------------------------------------
$ws2.ListObjects.Item(1).Resize($ws2.Range("A9:D59"))
$ws5.ListObjects.Item(1).Resize($ws5.Range("A9:D71"))

# ---------------------------------------------------------------------------
# Restore the on-screen selection for each sheet, then leave the originally
# active sheet (P3 - Ivar) selected so the workbook's active tab is unchanged.
# ---------------------------------------------------------------------------
$ws1.Activate()
$null = $ws1.Range("C58").Select()

$ws2.Activate()
$null = $ws2.Range("C62").Select()

$ws4.Activate()
$null = $ws4.Range("C58").Select()

$ws5.Activate()
$null = $ws5.Range("C77").Select()

$ws3.Activate()
$null = $ws3.Range("A82").Select()
